$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 6404
$ws.Range("I111").Value = 1609.3334
$ws.Range("K111").Value = 4828.0002
$ws.Range("M111").Value = -1761.0002
$ws.Range("H132").Value = 4524.4707
$ws.Range("I132").Value = 4927.3335
$ws.Range("K132").Value = 14782.0005
$ws.Range("M132").Value = -12252.0005
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140
$ws.Range("H137").Value = 1521.525
$ws.Range("I137").Value = 1467.7
$ws.Range("J137").Value = 1683
$ws.Range("K137").Value = 4403.1
$ws.Range("L137").Value = 5049
$ws.Range("M137").Value = -1853.1
$ws.Range("N137").Value = -10149

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2567.9058
$ws.Range("I32").Value = 1957.841
$ws.Range("J32").Value = 5550.4443
$ws.Range("K32").Value = 1957.841
$ws.Range("L32").Value = 5550.4443
$ws.Range("M32").Value = -1670.841
$ws.Range("N32").Value = -6124.4443
$ws.Range("H74").Value = 1812.9474
$ws.Range("I74").Value = 1885.0588
$ws.Range("K74").Value = 1885.0588
$ws.Range("M74").Value = -1011.0588
$ws.Range("H77").Value = 1812.9474
$ws.Range("I77").Value = 1885.0588
$ws.Range("K77").Value = 9425.294
$ws.Range("M77").Value = -5057.294
$ws.Range("H102").Value = 4131.125
$ws.Range("I102").Value = 2109.8
$ws.Range("J102").Value = 7500
$ws.Range("K102").Value = 2109.8
$ws.Range("L102").Value = 7500
$ws.Range("M102").Value = -487.8000000000002
$ws.Range("N102").Value = -10744
$ws.Range("H132").Value = 23487.791
$ws.Range("I132").Value = 2370.077
$ws.Range("K132").Value = 7110.231000000001
$ws.Range("M132").Value = -4580.231000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 360
$ws.Range("I12").Value = 360
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 360
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -192
$ws.Range("N12").ClearContents()
$ws.Range("H86").Value = 2017.619
$ws.Range("I86").Value = 1821.7646
$ws.Range("J86").Value = 2850
$ws.Range("K86").Value = 1821.7646
$ws.Range("L86").Value = 2850
$ws.Range("M86").Value = -698.7646
$ws.Range("N86").Value = -5096
$ws.Range("H89").Value = 2017.619
$ws.Range("I89").Value = 1821.7646
$ws.Range("J89").Value = 2850
$ws.Range("K89").Value = 9108.823
$ws.Range("L89").Value = 14250
$ws.Range("M89").Value = -3492.823
$ws.Range("N89").Value = -25482
$ws.Range("H94").Value = 3655.6956
$ws.Range("I94").Value = 1769.4286
$ws.Range("J94").Value = 6589.8887
$ws.Range("K94").Value = 1769.4286
$ws.Range("L94").Value = 6589.8887
$ws.Range("M94").Value = -1318.4286
$ws.Range("N94").Value = -7491.8887
$ws.Range("H99").Value = 2105.8
$ws.Range("I99").Value = 1815.9231
$ws.Range("K99").Value = 1815.9231
$ws.Range("M99").Value = -317.9231
$ws.Range("H105").Value = 1978.5883
$ws.Range("I105").Value = 1894.8
$ws.Range("K105").Value = 1894.8
$ws.Range("M105").Value = -147.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1100
$ws.Range("I16").Value = 1300
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1300
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1013
$ws.Range("N16").Value = -1574
$ws.Range("H105").Value = 12500816
$ws.Range("I105").Value = 15625517
$ws.Range("J105").Value = 2011
$ws.Range("K105").Value = 15625517
$ws.Range("L105").Value = 2011
$ws.Range("M105").Value = -15623770
$ws.Range("N105").Value = -5505
$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 1300
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1300
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 870
$ws.Range("N113").Value = -5340
$ws.Range("H134").Value = 1407.125
$ws.Range("I134").Value = 1187.5
$ws.Range("J134").Value = 1626.75
$ws.Range("K134").Value = 3562.5
$ws.Range("L134").Value = 4880.25
$ws.Range("M134").Value = -1027.5
$ws.Range("N134").Value = -9950.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 723.34
$ws.Range("J131").Value = 740.875
$ws.Range("L131").Value = 2222.625
$ws.Range("N131").Value = -12302.625

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6969230
$ws.Range("I11").Value = 9111111
$ws.Range("J11").Value = 2149997
$ws.Range("K11").Value = 9111111
$ws.Range("L11").Value = 2149997
$ws.Range("M11").Value = -9110972
$ws.Range("N11").Value = -2150275
$ws.Range("H80").Value = 3677.353
$ws.Range("I80").Value = 3236.25
$ws.Range("J80").Value = 4736
$ws.Range("K80").Value = 3236.25
$ws.Range("L80").Value = 4736
$ws.Range("M80").Value = -2238.25
$ws.Range("N80").Value = -6732
$ws.Range("H83").Value = 3677.353
$ws.Range("I83").Value = 3236.25
$ws.Range("J83").Value = 4736
$ws.Range("K83").Value = 16181.25
$ws.Range("L83").Value = 23680
$ws.Range("M83").Value = -11189.25
$ws.Range("N83").Value = -33664
$ws.Range("H93").Value = 15250
$ws.Range("J93").Value = 15250
$ws.Range("L93").Value = 15250
$ws.Range("N93").Value = -18994
$ws.Range("H103").Value = 32302
$ws.Range("J103").Value = 32302
$ws.Range("L103").Value = 32302
$ws.Range("N103").Value = -34646
$ws.Range("H113").Value = 2781
$ws.Range("I113").Value = 1813
$ws.Range("J113").Value = 4475
$ws.Range("K113").Value = 1813
$ws.Range("L113").Value = 4475
$ws.Range("M113").Value = 357
$ws.Range("N113").Value = -8815

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2328.5715
$ws.Range("I82").Value = 3368.75
$ws.Range("J82").Value = 941.6667
$ws.Range("K82").Value = 3368.75
$ws.Range("L82").Value = 941.6667
$ws.Range("M82").Value = -3007.75
$ws.Range("N82").Value = -1663.6667
$ws.Range("H85").Value = 2328.5715
$ws.Range("I85").Value = 3368.75
$ws.Range("J85").Value = 941.6667
$ws.Range("K85").Value = 3368.75
$ws.Range("L85").Value = 941.6667
$ws.Range("M85").Value = -2120.75
$ws.Range("N85").Value = -3437.6667
$ws.Range("H100").Value = 2637.0908
$ws.Range("I100").Value = 1791.3334
$ws.Range("K100").Value = 1791.3334
$ws.Range("M100").Value = -1250.3334
$ws.Range("H132").Value = 4832.6665
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4832.6665
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 14497.9995
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -19557.9995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 70007
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H132").Value = 3936.875
$ws.Range("J132").Value = 3927.8572
$ws.Range("L132").Value = 11783.5716
$ws.Range("N132").Value = -16843.5716
